$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")
$ws.Activate()

$ws.Range("B2").Value = "42"
$ws.Range("D2").Value = "Automation3"

$ws.Range("D2").Select()
